# Auto-update draw results: append the 2025-11-07 Pick 4 draw as a new
# row at the bottom of the Results sheet (row 52), keeping every value
# as plain text (matching the existing rows, which are all stored as
# text strings, not real numbers/dates).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 52
$rowRange = "A${newRow}:E${newRow}"

# Force text formatting first so Excel does not auto-coerce the
# date-looking / number-looking strings into real dates or numbers.
$ws.Range($rowRange).NumberFormat = "@"

$ws.Range("A${newRow}").Value = "2025-11-07"
$ws.Range("B${newRow}").Value = "Pick 4"
$ws.Range("C${newRow}").Value = "251107"
$ws.Range("D${newRow}").Value = "9-8-4-0"
$ws.Range("E${newRow}").Value = "2025-11-07T21:38:44.145+04:00"
